# Apply the changes described by the diff:
#  - Column Q (17th column) width 8 -> 10
#  - Column E (Statut) values "NA" -> "RAA" for rows 2..26
#  - Add a blank separator row (27) then 3 summary rows (28-30)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen column Q (index 17) from 8 to 10 ---
$ws.Columns.Item(17).ColumnWidth = 9.17

# --- Replace "NA" with "RAA" in the Statut column (E2:E26) ---
for ($r = 2; $r -le 26; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    if ($cell.Value2 -eq "NA") {
        $cell.Value = "RAA"
    }
}

# --- Append summary rows 27-30, reusing the formatting of the last data row ---
$lastRow = $ws.Range("A26:T26")
$newRows = $ws.Range("A27:T30")
$lastRow.Copy()
$newRows.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 27 stays blank (already got formats only, no values copied).

# Row 28 - Nombre de reçus
$ws.Range("A28").Value = "Nombre de reçus"
$ws.Range("B28").Value = 25

# Row 29 - Total avec reçus
$ws.Range("A29").Value = "Total avec reçus"
$ws.Range("B29").Value = 2902

# Row 30 - Total sans reçus
$ws.Range("A30").Value = "Total sans reçus"
$ws.Range("B30").Value = 0
